$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 86-87; everything currently at row 86 onward
# (through the former last row, 106) shifts down by two rows, landing at
# 88-108, carrying its original values with it unchanged.
$ws.Rows("86:87").Insert()

# Row 86: new weekly price entry (Murcott / Primera)
$ws.Range("A86").Value = 4
$ws.Range("B86").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C86").Value = "Los Lagos"
$ws.Range("D86").Value = 44474
$ws.Range("E86").Value = 10
$ws.Range("F86").Value = "Fruta"
$ws.Range("G86").Value = 100102
$ws.Range("H86").Value = "Cítricos"
$ws.Range("I86").Value = 100102004
$ws.Range("J86").Value = "Mandarina"
$ws.Range("K86").Value = "Murcott"
$ws.Range("L86").Value = "Primera"
$ws.Range("M86").Value = 600
$ws.Range("N86").Value = 6000
$ws.Range("O86").Value = 6500
$ws.Range("P86").Value = 6250
$ws.Range("Q86").Value = "$/bandeja 10 kilos"
$ws.Range("R86").Value = "Provincia de Limarí"
$ws.Range("S86").Value = 625
$ws.Range("T86").Value = 10

# Row 87: new weekly price entry (Murcott / Segunda)
$ws.Range("A87").Value = 4
$ws.Range("B87").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C87").Value = "Los Lagos"
$ws.Range("D87").Value = 44474
$ws.Range("E87").Value = 10
$ws.Range("F87").Value = "Fruta"
$ws.Range("G87").Value = 100102
$ws.Range("H87").Value = "Cítricos"
$ws.Range("I87").Value = 100102004
$ws.Range("J87").Value = "Mandarina"
$ws.Range("K87").Value = "Murcott"
$ws.Range("L87").Value = "Segunda"
$ws.Range("M87").Value = 200
$ws.Range("N87").Value = 4500
$ws.Range("O87").Value = 4500
$ws.Range("P87").Value = 4500
$ws.Range("Q87").Value = "$/bandeja 10 kilos"
$ws.Range("R87").Value = "Provincia de Limarí"
$ws.Range("S87").Value = 450
$ws.Range("T87").Value = 10
